$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be stored as literal text (avoids Excel auto-converting
# numeric-looking strings like "244.99" into numbers), while keeping the cell
# style identical to the untouched default style (no formatting side effects).
function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range($addr).Value = "'" + $val
    $sheet.Range($addr).Style = $sheet.Range("A1").Style
}

Set-TextValue $ws "D2" '36.924.70'
Set-TextValue $ws "E2" '  -0.50%  '

Set-TextValue $ws "D3" '2.037.66'
Set-TextValue $ws "E3" '  -0.65%  '

Set-TextValue $ws "E4" '  +0.20%  '

Set-TextValue $ws "D5" '244.99'
Set-TextValue $ws "E5" '  -1.58%  '

Set-TextValue $ws "E6" '  -0.90%  '

Set-TextValue $ws "D7" '58.41'
Set-TextValue $ws "E7" '  -1.40%  '

Set-TextValue $ws "E9" '  -1.74%  '

Set-TextValue $ws "D10" '0.0768'
Set-TextValue $ws "E10" '  -1.64%  '

Set-TextValue $ws "E11" '  +2.06%  '

Set-TextValue $ws "D12" '15.33'
Set-TextValue $ws "E12" '  -3.17%  '

Set-TextValue $ws "D13" '0.876'
Set-TextValue $ws "E13" '  +7.96%  '

Set-TextValue $ws "D14" '2.335.47'
Set-TextValue $ws "E14" '  -0.61%  '

Set-TextValue $ws "D15" '5.61'
Set-TextValue $ws "E15" '  +1.61%  '

Set-TextValue $ws "D16" '2.036.51'
Set-TextValue $ws "E16" '  -0.72%  '

Set-TextValue $ws "D17" '18.15'
Set-TextValue $ws "E17" '  +8.12%  '

Set-TextValue $ws "D18" '36.899.24'
Set-TextValue $ws "E18" '  -0.70%  '

Set-TextValue $ws "D19" '73.53'
Set-TextValue $ws "E19" '  -1.65%  '

Set-TextValue $ws "E20" '  -1.27%  '

Set-TextValue $ws "D21" '5.36'
Set-TextValue $ws "E21" '  +0.40%  '

Set-TextValue $ws "D22" '235.23'
Set-TextValue $ws "E22" '  -0.52%  '

Set-TextValue $ws "D23" '1.00'

Set-TextValue $ws "D24" '2.45'
Set-TextValue $ws "E24" '  +2.88%  '

Set-TextValue $ws "E25" '  +4.98%  '

Set-TextValue $ws "D26" '168.44'
Set-TextValue $ws "E26" '  +0.20%  '

Set-TextValue $ws "E27" '  -3.81%  '

Set-TextValue $ws "E28" '  +0.42%  '

Set-TextValue $ws "D29" '5.41'
Set-TextValue $ws "E29" '  +15.09%  '

Set-TextValue $ws "E31" '  -3.52%  '

Set-TextValue $ws "D32" '4.73'
Set-TextValue $ws "E32" '  +6.49%  '

Set-TextValue $ws "D33" '0.0611'
Set-TextValue $ws "E33" '  -0.24%  '

Set-TextValue $ws "E34" '  +0.19%  '

Set-TextValue $ws "D35" '0.0862'
Set-TextValue $ws "E35" '  -4.90%  '

Set-TextValue $ws "E36" '  +6.33%  '

Set-TextValue $ws "D37" '2.23'
Set-TextValue $ws "E37" '  +0.27%  '

Set-TextValue $ws "E38" '  -3.49%  '

Set-TextValue $ws "B39" 'THORChain'
Set-TextValue $ws "C39" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws "D39" '5.16'
Set-TextValue $ws "E39" '  +0.67%  '

Set-TextValue $ws "B40" 'HuobiToken'
Set-TextValue $ws "C40" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws "D40" '3.10'
Set-TextValue $ws "E40" '  -2.66%  '

Set-TextValue $ws "E41" '  +0.22%  '

Set-TextValue $ws "D42" '0.0960'
Set-TextValue $ws "E42" '  -11.63%  '

Set-TextValue $ws "E43" '  +0.70%  '

Set-TextValue $ws "D44" '96.99'
Set-TextValue $ws "E44" '  +1.16%  '

Set-TextValue $ws "D45" '16.85'
Set-TextValue $ws "E45" '  -3.68%  '

Set-TextValue $ws "D46" '1.291.25'
Set-TextValue $ws "E46" '  +1.04%  '

Set-TextValue $ws "D47" '2.33'
Set-TextValue $ws "E47" '  -4.06%  '

Set-TextValue $ws "B48" 'MXToken'
Set-TextValue $ws "C48" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws "D48" '2.85'
Set-TextValue $ws "E48" '  -0.19%  '

Set-TextValue $ws "B49" 'FTXToken'
Set-TextValue $ws "C49" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws "D49" '3.74'
Set-TextValue $ws "E49" '  +9.69%  '

Set-TextValue $ws "D50" '6.69'
Set-TextValue $ws "E50" '  +0.11%  '

Set-TextValue $ws "D51" '2.221.86'
Set-TextValue $ws "E51" '  -0.75%  '
